$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 11 with the new test-mail entry ---
$logs = $wb.Worksheets.Item("Logs")
$logs.Range("A11").Value = "Mijn product is kapot aangekomen. Kunnen jullie dit oplossen?"
$logs.Range("B11").Value = "mailmind.test@zohomail.eu"
$logs.Range("C11").Value = "Testmail #1: Mijn product is kapot aangekomen. Kunnen jullie dit oplossen?"
$logs.Range("D11").Value = "Retour / Terugbetaling"
$logs.Range("E11").Value = "Beste klant,`nBedankt voor uw bericht. Wat vervelend om te horen dat uw product kapot is aangekomen. Om dit probleem op te lossen, hebben we wat meer informatie nodig.`nZou u ons kunnen voorzien van uw ordernummer en eventueel foto's van het beschadigde product? Op die manier kunnen we uw situatie beter begrijpen en u verder helpen met een passende oplossing.`nAlvast bedankt voor uw medewerking.`nMet vriendelijke groet,`n[Naam] E-mailassistent"
$logs.Range("F11").Value = "2025-07-23 22:15:32"
$logs.Range("G11").Value = "Ja"
$logs.Range("H11").Value = "Nee"
$logs.Range("I11").Value = "Ja"
$logs.Range("J11").Value = "Nee"

# Multi-line content in E11 triggers Excel's automatic row-height growth;
# AutoFit the row back so it matches the sheet's other (un-customized) rows.
$logs.Rows.Item(11).AutoFit()

# --- Extend conditional formatting ranges on the Logs sheet to include row 11 ---
$cols = @("D","G","H","I","J")
foreach ($col in $cols) {
    $oldRange = $logs.Range($col + "2:" + $col + "10")
    $newRange = $logs.Range($col + "2:" + $col + "11")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: append aggregated row 3 for the new category ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A3").Value = "Retour / Terugbetaling"
$dash.Range("B3").Value = 1

# --- Update the chart series so it covers the new Dashboard rows (A2:A3 / B2:B3) ---
$co = $dash.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$3,Dashboard!`$B`$2:`$B`$3,1)"
